$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6250
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 6250
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 6250
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -13134
$ws.Range("H125").Value = 2244
$ws.Range("I125").Value = 2021.3334
$ws.Range("J125").Value = 2466.6667
$ws.Range("K125").Value = 18192.0006
$ws.Range("L125").Value = 22200.0003
$ws.Range("M125").Value = -15732.0006
$ws.Range("N125").Value = -27120.0003
$ws.Range("H132").Value = 20723.348
$ws.Range("I132").Value = 2853.0264
$ws.Range("K132").Value = 8559.0792
$ws.Range("M132").Value = -6029.0792
$ws.Range("H137").Value = 4865
$ws.Range("I137").Value = 990.5333000000001
$ws.Range("J137").Value = 9335.538
$ws.Range("K137").Value = 2971.5999
$ws.Range("L137").Value = 28006.614
$ws.Range("M137").Value = -421.5999000000002
$ws.Range("N137").Value = -33106.614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1405.56
$ws.Range("I2").Value = 1497.1364
$ws.Range("J2").Value = 734
$ws.Range("K2").Value = 1497.1364
$ws.Range("L2").Value = 734
$ws.Range("M2").Value = -1384.1364
$ws.Range("N2").Value = -960
$ws.Range("H53").Value = 12779.75
$ws.Range("I53").Value = 10891.143
$ws.Range("J53").Value = 26000
$ws.Range("K53").Value = 10891.143
$ws.Range("L53").Value = 26000
$ws.Range("M53").Value = -10209.143
$ws.Range("N53").Value = -27364
$ws.Range("H74").Value = 2494.8518
$ws.Range("I74").Value = 2094.1304
$ws.Range("K74").Value = 2094.1304
$ws.Range("M74").Value = -1220.1304
$ws.Range("H77").Value = 2494.8518
$ws.Range("I77").Value = 2094.1304
$ws.Range("K77").Value = 10470.652
$ws.Range("M77").Value = -6102.652
$ws.Range("H116").Value = 1405.56
$ws.Range("I116").Value = 1497.1364
$ws.Range("J116").Value = 734
$ws.Range("K116").Value = 1497.1364
$ws.Range("L116").Value = 734
$ws.Range("M116").Value = 796.8635999999999
$ws.Range("N116").Value = -5322
$ws.Range("H122").Value = 1372.3158
$ws.Range("I122").Value = 1460.92
$ws.Range("J122").Value = 1201.9231
$ws.Range("K122").Value = 4382.76
$ws.Range("L122").Value = 3605.7693
$ws.Range("M122").Value = -1932.76
$ws.Range("N122").Value = -8505.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1405.56
$ws.Range("I3").Value = 1497.1364
$ws.Range("J3").Value = 734
$ws.Range("K3").Value = 1497.1364
$ws.Range("L3").Value = 734
$ws.Range("M3").Value = -1383.1364
$ws.Range("N3").Value = -962
$ws.Range("H64").Value = 603.9286
$ws.Range("I64").Value = 699.4286
$ws.Range("J64").Value = 508.42856
$ws.Range("K64").Value = 699.4286
$ws.Range("L64").Value = 508.42856
$ws.Range("M64").Value = -474.4286
$ws.Range("N64").Value = -958.4285600000001
$ws.Range("H67").Value = 603.9286
$ws.Range("I67").Value = 699.4286
$ws.Range("J67").Value = 508.42856
$ws.Range("K67").Value = 699.4286
$ws.Range("L67").Value = 508.42856
$ws.Range("M67").Value = 80.57140000000004
$ws.Range("N67").Value = -2068.42856
$ws.Range("H99").Value = 1982.9117
$ws.Range("I99").Value = 1524.76
$ws.Range("J99").Value = 3255.5557
$ws.Range("K99").Value = 1524.76
$ws.Range("L99").Value = 3255.5557
$ws.Range("M99").Value = -26.75999999999999
$ws.Range("N99").Value = -6251.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4163.8877
$ws.Range("I31").Value = 2090.5833
$ws.Range("J31").Value = 4836.311
$ws.Range("K31").Value = 2090.5833
$ws.Range("L31").Value = 4836.311
$ws.Range("M31").Value = -1795.5833
$ws.Range("N31").Value = -5426.311
$ws.Range("H34").Value = 4163.8877
$ws.Range("I34").Value = 2090.5833
$ws.Range("J34").Value = 4836.311
$ws.Range("K34").Value = 2090.5833
$ws.Range("L34").Value = 4836.311
$ws.Range("M34").Value = -1888.5833
$ws.Range("N34").Value = -5240.311
$ws.Range("H122").Value = 150739.25
$ws.Range("I122").Value = 150739.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 452217.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -449767.75
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 315.5484
$ws.Range("I12").Value = 244.14285
$ws.Range("J12").Value = 336.375
$ws.Range("K12").Value = 732.4285500000001
$ws.Range("L12").Value = 1009.125
$ws.Range("M12").Value = -559.4285500000001
$ws.Range("N12").Value = -1355.125
$ws.Range("H29").Value = 7651
$ws.Range("I29").Value = 537.1429000000001
$ws.Range("J29").Value = 14764.857
$ws.Range("K29").Value = 1611.4287
$ws.Range("L29").Value = 44294.571
$ws.Range("M29").Value = -1334.4287
$ws.Range("N29").Value = -44848.571
$ws.Range("H68").Value = 4461.375
$ws.Range("J68").Value = 5717.1665
$ws.Range("L68").Value = 17151.4995
$ws.Range("N68").Value = -18773.4995
$ws.Range("H71").Value = 4461.375
$ws.Range("J71").Value = 5717.1665
$ws.Range("L71").Value = 51454.4985
$ws.Range("N71").Value = -59566.4985
$ws.Range("H113").Value = 3140.9
$ws.Range("I113").Value = 4875.2173
$ws.Range("K113").Value = 14625.6519
$ws.Range("M113").Value = -12455.6519
$ws.Range("H132").Value = 5376.143
$ws.Range("I132").Value = 2017.4
$ws.Range("K132").Value = 18156.6
$ws.Range("M132").Value = -15626.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4716.7144
$ws.Range("I43").Value = 3017
$ws.Range("K43").Value = 3017
$ws.Range("M43").Value = -2866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3454.3076
$ws.Range("I7").Value = 2414.5715
$ws.Range("K7").Value = 2414.5715
$ws.Range("M7").Value = -2302.5715
$ws.Range("H40").Value = 4580.381
$ws.Range("I40").Value = 3147.2
$ws.Range("J40").Value = 8163.3335
$ws.Range("K40").Value = 3147.2
$ws.Range("L40").Value = 8163.3335
$ws.Range("M40").Value = -3011.2
$ws.Range("N40").Value = -8435.333500000001
$ws.Range("H55").Value = 746.91895
$ws.Range("I55").Value = 770.6
$ws.Range("J55").Value = 730.7727
$ws.Range("K55").Value = 770.6
$ws.Range("L55").Value = 730.7727
$ws.Range("M55").Value = -597.6
$ws.Range("N55").Value = -1076.7727
$ws.Range("H100").Value = 2581.6667
$ws.Range("I100").Value = 2538
$ws.Range("K100").Value = 2538
$ws.Range("M100").Value = -1997
$ws.Range("H126").Value = 3454.3076
$ws.Range("I126").Value = 2414.5715
$ws.Range("K126").Value = 7243.7145
$ws.Range("M126").Value = -4773.7145
$ws.Range("H136").Value = 1992.4286
$ws.Range("I136").Value = 1513.0454
$ws.Range("J136").Value = 3750.1667
$ws.Range("K136").Value = 4539.1362
$ws.Range("L136").Value = 11250.5001
$ws.Range("M136").Value = -1989.1362
$ws.Range("N136").Value = -16350.5001
$ws.Range("H137").Value = 47750
$ws.Range("J137").Value = 47750
$ws.Range("L137").Value = 47750
$ws.Range("N137").Value = -57950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1500
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = ""
$ws.Range("N96").Value = -4246
